$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 4015
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(72, 8).Value = 4015
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(125, 8).Value = 1209.0769
$ws.Cells.Item(125, 9).Value = 1117.7142
$ws.Cells.Item(125, 10).Value = 1315.6666
$ws.Cells.Item(125, 11).Value = 10059.4278
$ws.Cells.Item(125, 12).Value = 11840.9994
$ws.Cells.Item(125, 13).Value = -7599.427799999999
$ws.Cells.Item(125, 14).Value = -16760.9994
$ws.Cells.Item(137, 8).Value = 5933.4717
$ws.Cells.Item(137, 9).Value = 4291.5527
$ws.Cells.Item(137, 10).Value = 10093
$ws.Cells.Item(137, 11).Value = 12874.6581
$ws.Cells.Item(137, 12).Value = 30279
$ws.Cells.Item(137, 13).Value = -10324.6581
$ws.Cells.Item(137, 14).Value = -35379
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13779.423
$ws.Cells.Item(32, 9).Value = 12880.178
$ws.Cells.Item(32, 10).Value = 19560.285
$ws.Cells.Item(32, 11).Value = 12880.178
$ws.Cells.Item(32, 12).Value = 19560.285
$ws.Cells.Item(32, 13).Value = -12593.178
$ws.Cells.Item(32, 14).Value = -20134.285
$ws.Cells.Item(36, 8).Value = 1000
$ws.Cells.Item(36, 9).Value = 1000
$ws.Cells.Item(36, 11).Value = 1000
$ws.Cells.Item(36, 13).Value = -654
$ws.Cells.Item(37, 8).Value = 37710
$ws.Cells.Item(37, 10).Value = 37710
$ws.Cells.Item(37, 12).Value = 37710
$ws.Cells.Item(37, 14).Value = -38256
$ws.Cells.Item(61, 8).Value = 1976.4706
$ws.Cells.Item(61, 9).Value = 1243.9474
$ws.Cells.Item(61, 10).Value = 2904.3333
$ws.Cells.Item(61, 11).Value = 1243.9474
$ws.Cells.Item(61, 12).Value = 2904.3333
$ws.Cells.Item(61, 13).Value = -1031.9474
$ws.Cells.Item(61, 14).Value = -3328.3333
$ws.Cells.Item(111, 8).Value = 49800
$ws.Cells.Item(111, 10).Value = 49800
$ws.Cells.Item(111, 12).Value = 49800
$ws.Cells.Item(111, 14).Value = -57980
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(125, 8).Value = 60000
$ws.Cells.Item(125, 10).Value = 60000
$ws.Cells.Item(125, 12).Value = 60000
$ws.Cells.Item(125, 14).Value = -69840
$ws.Cells.Item(132, 8).Value = 29415516
$ws.Cells.Item(132, 9).Value = 55559600
$ws.Cells.Item(132, 10).Value = 3423
$ws.Cells.Item(132, 11).Value = 166678800
$ws.Cells.Item(132, 12).Value = 10269
$ws.Cells.Item(132, 13).Value = -166676270
$ws.Cells.Item(132, 14).Value = -15329
$ws.Cells.Item(136, 8).Value = 1976.4706
$ws.Cells.Item(136, 9).Value = 1243.9474
$ws.Cells.Item(136, 10).Value = 2904.3333
$ws.Cells.Item(136, 11).Value = 3731.8422
$ws.Cells.Item(136, 12).Value = 8712.999899999999
$ws.Cells.Item(136, 13).Value = -1181.8422
$ws.Cells.Item(136, 14).Value = -13812.9999
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(103, 8).Value = 200000
$ws.Cells.Item(103, 10).Value = 200000
$ws.Cells.Item(103, 12).Value = 200000
$ws.Cells.Item(103, 14).Value = -202344
$ws.Cells.Item(105, 8).Value = 2559.6
$ws.Cells.Item(105, 9).Value = 2538.3
$ws.Cells.Item(105, 10).Value = 2602.2
$ws.Cells.Item(105, 11).Value = 2538.3
$ws.Cells.Item(105, 12).Value = 2602.2
$ws.Cells.Item(105, 13).Value = -791.3000000000002
$ws.Cells.Item(105, 14).Value = -6096.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1635.1666
$ws.Cells.Item(16, 9).Value = 2002.75
$ws.Cells.Item(16, 10).Value = 900
$ws.Cells.Item(16, 11).Value = 2002.75
$ws.Cells.Item(16, 12).Value = 900
$ws.Cells.Item(16, 13).Value = -1715.75
$ws.Cells.Item(16, 14).Value = -1474
$ws.Cells.Item(31, 8).Value = 4170555
$ws.Cells.Item(31, 9).Value = 2188.7036
$ws.Cells.Item(31, 11).Value = 2188.7036
$ws.Cells.Item(31, 13).Value = -1893.7036
$ws.Cells.Item(34, 8).Value = 4170555
$ws.Cells.Item(34, 9).Value = 2188.7036
$ws.Cells.Item(34, 11).Value = 2188.7036
$ws.Cells.Item(34, 13).Value = -1986.7036
$ws.Cells.Item(36, 8).Value = 35016
$ws.Cells.Item(36, 9).Value = 50774
$ws.Cells.Item(36, 10).Value = 3500
$ws.Cells.Item(36, 11).Value = 50774
$ws.Cells.Item(36, 12).Value = 3500
$ws.Cells.Item(36, 13).Value = -50386
$ws.Cells.Item(36, 14).Value = -4276
$ws.Cells.Item(40, 8).Value = 35016
$ws.Cells.Item(40, 9).Value = 50774
$ws.Cells.Item(40, 10).Value = 3500
$ws.Cells.Item(40, 11).Value = 50774
$ws.Cells.Item(40, 12).Value = 3500
$ws.Cells.Item(40, 13).Value = -50614
$ws.Cells.Item(40, 14).Value = -3820
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 552.5946
$ws.Cells.Item(107, 9).Value = 361.56522
$ws.Cells.Item(107, 10).Value = 866.4286
$ws.Cells.Item(107, 11).Value = 361.56522
$ws.Cells.Item(107, 12).Value = 866.4286
$ws.Cells.Item(107, 13).Value = 1558.43478
$ws.Cells.Item(107, 14).Value = -4706.4286
$ws.Cells.Item(113, 8).Value = 1635.1666
$ws.Cells.Item(113, 9).Value = 2002.75
$ws.Cells.Item(113, 10).Value = 900
$ws.Cells.Item(113, 11).Value = 2002.75
$ws.Cells.Item(113, 12).Value = 900
$ws.Cells.Item(113, 13).Value = 167.25
$ws.Cells.Item(113, 14).Value = -5240
$ws.Cells.Item(122, 8).Value = 58190.094
$ws.Cells.Item(122, 9).Value = 71481.06
$ws.Cells.Item(122, 10).Value = 1703.5
$ws.Cells.Item(122, 11).Value = 214443.18
$ws.Cells.Item(122, 12).Value = 5110.5
$ws.Cells.Item(122, 13).Value = -211993.18
$ws.Cells.Item(122, 14).Value = -10010.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1288.2222
$ws.Cells.Item(34, 9).Value = 320.6
$ws.Cells.Item(34, 10).Value = 1660.3846
$ws.Cells.Item(34, 11).Value = 961.8000000000001
$ws.Cells.Item(34, 12).Value = 4981.1538
$ws.Cells.Item(34, 13).Value = -877.8000000000001
$ws.Cells.Item(34, 14).Value = -5149.1538
$ws.Cells.Item(35, 8).Value = 1500
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 1500
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 4500
$ws.Cells.Item(35, 13).ClearContents()
$ws.Cells.Item(35, 14).Value = -5076
$ws.Cells.Item(46, 8).Value = 1065.7
$ws.Cells.Item(46, 9).Value = 67.666664
$ws.Cells.Item(46, 10).Value = 1493.4286
$ws.Cells.Item(46, 11).Value = 202.999992
$ws.Cells.Item(46, 12).Value = 4480.2858
$ws.Cells.Item(46, 13).Value = -111.999992
$ws.Cells.Item(46, 14).Value = -4662.2858
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5227.433
$ws.Cells.Item(80, 9).Value = 4112.0557
$ws.Cells.Item(80, 10).Value = 6900.5
$ws.Cells.Item(80, 11).Value = 4112.0557
$ws.Cells.Item(80, 12).Value = 6900.5
$ws.Cells.Item(80, 13).Value = -3114.0557
$ws.Cells.Item(80, 14).Value = -8896.5
$ws.Cells.Item(83, 8).Value = 5227.433
$ws.Cells.Item(83, 9).Value = 4112.0557
$ws.Cells.Item(83, 10).Value = 6900.5
$ws.Cells.Item(83, 11).Value = 20560.2785
$ws.Cells.Item(83, 12).Value = 34502.5
$ws.Cells.Item(83, 13).Value = -15568.2785
$ws.Cells.Item(83, 14).Value = -44486.5
$ws.Cells.Item(132, 8).Value = 3140.087
$ws.Cells.Item(132, 9).Value = 2162.2
$ws.Cells.Item(132, 10).Value = 3892.3076
$ws.Cells.Item(132, 11).Value = 6486.599999999999
$ws.Cells.Item(132, 12).Value = 11676.9228
$ws.Cells.Item(132, 13).Value = -3956.599999999999
$ws.Cells.Item(132, 14).Value = -16736.9228
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 14).ClearContents()
$ws.Cells.Item(7, 8).Value = 111115000
$ws.Cells.Item(7, 9).Value = 500001120
$ws.Cells.Item(7, 11).Value = 500001120
$ws.Cells.Item(7, 13).Value = -500001008
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).ClearContents()
$ws.Cells.Item(41, 8).Value = 28995
$ws.Cells.Item(41, 10).Value = 28995
$ws.Cells.Item(41, 12).Value = 28995
$ws.Cells.Item(41, 14).Value = -29871
$ws.Cells.Item(45, 8).Value = 9672.5
$ws.Cells.Item(45, 9).Value = 7000
$ws.Cells.Item(45, 10).Value = 10563.333
$ws.Cells.Item(45, 11).Value = 7000
$ws.Cells.Item(45, 12).Value = 10563.333
$ws.Cells.Item(45, 13).Value = -6593
$ws.Cells.Item(45, 14).Value = -11377.333
$ws.Cells.Item(68, 8).Value = 3946.6667
$ws.Cells.Item(68, 9).Value = 4000
$ws.Cells.Item(68, 10).Value = 3900
$ws.Cells.Item(68, 11).Value = 4000
$ws.Cells.Item(68, 12).Value = 3900
$ws.Cells.Item(68, 13).Value = -3251
$ws.Cells.Item(68, 14).Value = -5398
$ws.Cells.Item(71, 8).Value = 3946.6667
$ws.Cells.Item(71, 9).Value = 4000
$ws.Cells.Item(71, 10).Value = 3900
$ws.Cells.Item(71, 11).Value = 20000
$ws.Cells.Item(71, 12).Value = 19500
$ws.Cells.Item(71, 13).Value = -16256
$ws.Cells.Item(71, 14).Value = -26988
$ws.Cells.Item(118, 8).Value = 28762.5
$ws.Cells.Item(118, 10).Value = 28762.5
$ws.Cells.Item(118, 12).Value = 28762.5
$ws.Cells.Item(118, 14).Value = -32076.5
$ws.Cells.Item(123, 8).Value = 60000
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 60000
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 60000
$ws.Cells.Item(123, 13).ClearContents()
$ws.Cells.Item(123, 14).Value = -69800
$ws.Cells.Item(126, 8).Value = 111115000
$ws.Cells.Item(126, 9).Value = 500001120
$ws.Cells.Item(126, 11).Value = 1500003360
$ws.Cells.Item(126, 13).Value = -1500000890
$ws.Cells.Item(132, 8).Value = 3533.0908
$ws.Cells.Item(132, 10).Value = 4274.0835
$ws.Cells.Item(132, 12).Value = 12822.2505
$ws.Cells.Item(132, 14).Value = -17882.2505
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 49420
$ws.Cells.Item(16, 10).Value = 49420
$ws.Cells.Item(16, 12).Value = 49420
$ws.Cells.Item(16, 14).Value = -50004
$ws.Cells.Item(28, 8).Value = 40000
$ws.Cells.Item(28, 10).Value = 40000
$ws.Cells.Item(28, 12).Value = 40000
$ws.Cells.Item(28, 14).Value = -40696
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()
$ws.Cells.Item(116, 8).Value = 37500
$ws.Cells.Item(116, 10).Value = 37500
$ws.Cells.Item(116, 12).Value = 37500
$ws.Cells.Item(116, 14).Value = -46678
$ws.Cells.Item(126, 8).Value = 2942446.8
$ws.Cells.Item(126, 9).Value = 2942446.8
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8827340.399999999
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -8824870.399999999
$ws.Cells.Item(126, 14).ClearContents()
